$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, pushing the existing rows 39-137 down to 40-138.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly record.
$ws.Range("A39").Value = 8
$ws.Range("B39").Value = "Terminal La Palmera de La Serena"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 44715
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 100112001
$ws.Range("G39").Value = "Berenjena"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = 8500
$ws.Range("N39").Value = "$/caja 50 unidades"
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 170
$ws.Range("Q39").Value = 50
$ws.Range("R39").Value = "Hortaliza"
